# chore: update Sheets via scheduled runner
# Refreshes the market-price-derived leve profit figures (columns H-N) for
# the affected leve rows across the per-job sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 708.5454999999999
$ws.Range("I4").Value = 421.33334
$ws.Range("J4").Value = 2001
$ws.Range("K4").Value = 421.33334
$ws.Range("L4").Value = 2001
$ws.Range("M4").Value = -307.33334
$ws.Range("N4").Value = -2229
$ws.Range("H33").Value = 301.33334
$ws.Range("I33").Value = 228
$ws.Range("K33").Value = 228
$ws.Range("M33").Value = 1
$ws.Range("H40").Value = 3183.5862
$ws.Range("I40").Value = 2779.111
$ws.Range("J40").Value = 3845.4546
$ws.Range("K40").Value = 2779.111
$ws.Range("L40").Value = 3845.4546
$ws.Range("M40").Value = -2604.111
$ws.Range("N40").Value = -4195.4546

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 721.7742
$ws.Range("I2").Value = 584.2963
$ws.Range("K2").Value = 584.2963
$ws.Range("M2").Value = -471.2963
$ws.Range("H22").Value = 3113.7144
$ws.Range("I22").Value = 8498
$ws.Range("J22").Value = 960
$ws.Range("K22").Value = 8498
$ws.Range("L22").Value = 960
$ws.Range("M22").Value = -8199
$ws.Range("N22").Value = -1558
$ws.Range("H32").Value = 3409.4023
$ws.Range("I32").Value = 2773.2026
$ws.Range("K32").Value = 2773.2026
$ws.Range("M32").Value = -2486.2026
$ws.Range("H74").Value = 2256.9092
$ws.Range("J74").Value = 3501
$ws.Range("L74").Value = 3501
$ws.Range("N74").Value = -5249
$ws.Range("H77").Value = 2256.9092
$ws.Range("J77").Value = 3501
$ws.Range("L77").Value = 17505
$ws.Range("N77").Value = -26241
$ws.Range("H97").Value = 1979.6364
$ws.Range("J97").Value = 6948.5713
$ws.Range("L97").Value = 6948.5713
$ws.Range("N97").Value = -7940.5713
$ws.Range("H116").Value = 721.7742
$ws.Range("I116").Value = 584.2963
$ws.Range("K116").Value = 584.2963
$ws.Range("M116").Value = 1709.7037
$ws.Range("H132").Value = 3668.2778
$ws.Range("I132").Value = 3521.6
$ws.Range("J132").Value = 4401.6665
$ws.Range("K132").Value = 10564.8
$ws.Range("L132").Value = 13204.9995
$ws.Range("M132").Value = -8034.799999999999
$ws.Range("N132").Value = -18264.9995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 721.7742
$ws.Range("I3").Value = 584.2963
$ws.Range("K3").Value = 584.2963
$ws.Range("M3").Value = -470.2963
$ws.Range("H20").Value = 1224.8572
$ws.Range("I20").Value = 1225
$ws.Range("J20").Value = 1224.6666
$ws.Range("K20").Value = 1225
$ws.Range("L20").Value = 1224.6666
$ws.Range("M20").Value = -978
$ws.Range("N20").Value = -1718.6666
$ws.Range("H54").Value = 3696.4
$ws.Range("I54").Value = 3696.4
$ws.Range("K54").Value = 3696.4
$ws.Range("M54").Value = -3212.4
$ws.Range("H64").Value = 1545.8572
$ws.Range("I64").Value = 1150.25
$ws.Range("K64").Value = 1150.25
$ws.Range("M64").Value = -925.25
$ws.Range("H67").Value = 1545.8572
$ws.Range("I67").Value = 1150.25
$ws.Range("K67").Value = 1150.25
$ws.Range("M67").Value = -370.25
$ws.Range("H86").Value = 5484.7036
$ws.Range("I86").Value = 5373.7
$ws.Range("J86").Value = 5801.857
$ws.Range("K86").Value = 5373.7
$ws.Range("L86").Value = 5801.857
$ws.Range("M86").Value = -4250.7
$ws.Range("N86").Value = -8047.857
$ws.Range("H89").Value = 5484.7036
$ws.Range("I89").Value = 5373.7
$ws.Range("J89").Value = 5801.857
$ws.Range("K89").Value = 26868.5
$ws.Range("L89").Value = 29009.285
$ws.Range("M89").Value = -21252.5
$ws.Range("N89").Value = -40241.285

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2431.2144
$ws.Range("I94").Value = 2231.5
$ws.Range("K94").Value = 2231.5
$ws.Range("M94").Value = -1780.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 17777.5
$ws.Range("I109").Value = 5555
$ws.Range("J109").Value = 30000
$ws.Range("K109").Value = 16665
$ws.Range("L109").Value = 90000
$ws.Range("M109").Value = -15625
$ws.Range("N109").Value = -92080

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3518.0908
$ws.Range("I80").Value = 3232.5715
$ws.Range("K80").Value = 3232.5715
$ws.Range("M80").Value = -2234.5715
$ws.Range("H83").Value = 3518.0908
$ws.Range("I83").Value = 3232.5715
$ws.Range("K83").Value = 16162.8575
$ws.Range("M83").Value = -11170.8575
$ws.Range("H132").Value = 3223
$ws.Range("I132").Value = 2511.1365
$ws.Range("K132").Value = 7533.4095
$ws.Range("M132").Value = -5003.4095

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 531
$ws.Range("I16").Value = 446.92
$ws.Range("K16").Value = 446.92
$ws.Range("M16").Value = -276.92
$ws.Range("H22").Value = 1488.5555
$ws.Range("J22").Value = 1665
$ws.Range("L22").Value = 1665
$ws.Range("N22").Value = -2255
$ws.Range("H27").Value = 1488.5555
$ws.Range("J27").Value = 1665
$ws.Range("L27").Value = 1665
$ws.Range("N27").Value = -1879
$ws.Range("H40").Value = 2880.6155
$ws.Range("I40").Value = 2886.182
$ws.Range("K40").Value = 2886.182
$ws.Range("M40").Value = -2750.182
$ws.Range("H93").Value = 2570
$ws.Range("I93").Value = 2373.5
$ws.Range("K93").Value = 2373.5
$ws.Range("M93").Value = -1125.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 51188
$ws.Range("J46").Value = 51188
$ws.Range("L46").Value = 51188
$ws.Range("N46").Value = -51650
$ws.Range("H62").Value = 9019.200000000001
$ws.Range("I62").Value = 9065.666999999999
$ws.Range("K62").Value = 9065.666999999999
$ws.Range("M62").Value = -8441.666999999999
$ws.Range("H65").Value = 9019.200000000001
$ws.Range("I65").Value = 9065.666999999999
$ws.Range("K65").Value = 45328.335
$ws.Range("M65").Value = -42208.335
$ws.Range("H122").Value = 31806.684
$ws.Range("I122").Value = 27517.533
$ws.Range("J122").Value = 47891
$ws.Range("K122").Value = 82552.599
$ws.Range("L122").Value = 143673
$ws.Range("M122").Value = -80102.599
$ws.Range("N122").Value = -148573
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530
$ws.Range("H134").Value = 51188
$ws.Range("J134").Value = 51188
$ws.Range("L134").Value = 153564
$ws.Range("N134").Value = -158634
